# Insert a new row at row 16 (shifts existing rows 16..115 down to 17..116)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new record's data
$ws.Range("A16").Value = 5
$ws.Range("B16").Value = "Macroferia Regional de Talca"
$ws.Range("C16").Value = "Maule"
$ws.Range("D16").Value = 44613
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = 100112030
$ws.Range("G16").Value = "Poroto granado"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 300
$ws.Range("K16").Value = 18000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 18000
$ws.Range("N16").Value = "$/saco 25 kilos"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 720
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = "Hortaliza"

# Ensure the date cell keeps the date-formatted style (same as other rows' column D)
$ws.Range("D16").NumberFormat = $ws.Range("D17").NumberFormat
